$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 14: 280-XX-22-A -> 280-011-22-A, fill in B/C, bump D to 3 ---
$ws.Range("A14").Value2 = "280-011-22-A"
$ws.Range("B14").Value2 = 1
$ws.Range("C14").Value2 = 0.095
$ws.Range("D14").Value2 = 3

# --- Insert three new rows before the old row 15 (BAJA-035) so it ---
# --- and the old row 16 (BAJA-049) shift down to rows 18 and 19.   ---
$ws.Rows("15:17").Insert()

# New row 15: 280-013-22-A
$ws.Range("A15").Value2 = "280-013-22-A"
$ws.Range("B15").Value2 = 1
$ws.Range("C15").Value2 = 0.065
$ws.Range("D15").Value2 = 3

# New row 16: 280-014-22-A
$ws.Range("A16").Value2 = "280-014-22-A"
$ws.Range("B16").Value2 = 1
$ws.Range("C16").Value2 = 0.065
$ws.Range("D16").Value2 = 3

# New row 17: 280-015-22-A
$ws.Range("A17").Value2 = "280-015-22-A"
$ws.Range("B17").Value2 = 1
$ws.Range("C17").Value2 = 0.065
$ws.Range("D17").Value2 = 6

# Keep the cursor where the author left it.
$ws.Range("D22").Select()
